$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

$ws.Range("E5").Value = "Test 1"
$ws.Range("F5").Value = "Test 2"
$ws.Range("E17").Value = "Test 1"
$ws.Range("F17").Value = "Test 2"
